$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 24: "git clone https://github.com/tsduck/tsduck"
#   -> split into "git clone https://" + "github.com/tsduck/tsduck.git"
# -----------------------------------------------------------------
$s24 = $p.Slides.Item(24)
$tr24 = $s24.Shapes.Item(1).TextFrame.TextRange
$para24 = $tr24.Paragraphs(3)
$prefixLen = 18  # "git clone https://"
$visibleLen = $para24.Text.Length - 1
$tail24 = $para24.Characters($prefixLen + 1, $visibleLen - $prefixLen)
$tail24.Text = "github.com/tsduck/tsduck.git"

# -----------------------------------------------------------------
# Slide 25: "Try to find a solution using existing TSDuck"
#   -> split into "Try to find a solution using existing " + "TSDuck"
# -----------------------------------------------------------------
$s25 = $p.Slides.Item(25)
$tr25 = $s25.Shapes.Item(1).TextFrame.TextRange
$para25 = $tr25.Paragraphs(2)
$prefixLen = 38  # "Try to find a solution using existing "
$visibleLen = $para25.Text.Length - 1
$tail25 = $para25.Characters($prefixLen + 1, $visibleLen - $prefixLen)
$tail25.Text = "TSDuck"

# -----------------------------------------------------------------
# Slide 27: "Known issue with DVB tuners & DirectShow on 64 bits"
#   -> split into "Known issue with DVB tuners & DirectShow on 64 " + "bits Windows"
# -----------------------------------------------------------------
$s27 = $p.Slides.Item(27)
$tr27 = $s27.Shapes.Item(1).TextFrame.TextRange
$para27 = $tr27.Paragraphs(7)
$prefixLen = 47  # "Known issue with DVB tuners & DirectShow on 64 "
$visibleLen = $para27.Text.Length - 1
$tail27 = $para27.Characters($prefixLen + 1, $visibleLen - $prefixLen)
$tail27.Text = "bits Windows"

# -----------------------------------------------------------------
# Slide 30: "... " + "application-specific rules ..." (2 runs)
#   -> merge into single run "... application-specific rules ..."
# -----------------------------------------------------------------
$s30 = $p.Slides.Item(30)
$tr30 = $s30.Shapes.Item(1).TextFrame.TextRange
$para30 = $tr30.Paragraphs(7)
$run1Len = 4  # "... "
$full30 = "... application-specific rules ..."
$run1_30 = $para30.Characters(1, $run1Len)
$run1_30.Text = $full30
$visibleLen = $para30.Text.Length - 1
$rest30 = $para30.Characters($full30.Length + 1, $visibleLen - $full30.Length)
$rest30.Text = ""

# -----------------------------------------------------------------
# Slide 6: "User's " + "Guide" (2 runs) -> merge into "User's Guide"
# -----------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tr6 = $s6.Shapes.Item(1).TextFrame.TextRange
$paraGuide = $tr6.Paragraphs(1)
$run1Len = 7  # "User's "
$fullGuide = "User's Guide"
$run1_guide = $paraGuide.Characters(1, $run1Len)
$run1_guide.Text = $fullGuide
$visibleLen = $paraGuide.Text.Length - 1
$restGuide = $paraGuide.Characters($fullGuide.Length + 1, $visibleLen - $fullGuide.Length)
$restGuide.Text = ""

# -----------------------------------------------------------------
# Slide 6: "generated " + "by Doxygen from source " + "code" (3 runs + endParaRPr)
#   -> merge into single run "generated by Doxygen from source code", endParaRPr removed
# -----------------------------------------------------------------
$paraDox = $tr6.Paragraphs(8)
$paraDox.Delete()
$paraNext = $tr6.Paragraphs(8)  # now "C++ common code reference"
$paraNext.InsertBefore("generated by Doxygen from source code" + [char]13)
